$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  ,('ECs', 'Efnb2', 'Ephb1', 'ECs', 3, 1, 37.44676833333333, 112.340305, 0.6029245723174422, 0.6029245723174423, 3, 1, 1.397474666666667, 4.192424, 0.9253320219903044, 0.9253320219903044, 52.33091009436888, 470.97819084932, 0.5579054136101383, 0.5579054136101385)
  ,('ECs', 'Efnb2', 'Ephb1', 'MuSCs', 3, 1, 37.44676833333333, 112.340305, 0.6029245723174422, 0.6029245723174423, 3, 1, 0.07803566666666666, 0.234107, 0.05167099121464913, 0.05167099121464913, 2.922183531403888, 26.299651782635, 0.03115371027931064, 0.03115371027931065)
  ,('ECs', 'Efnb2', 'Ephb1', 'Neutrophils', 3, 1, 37.44676833333333, 112.340305, 0.6029245723174422, 0.6029245723174423, 1, 0.3333333333333333, 0.034731, 0.104193, 0.02299698679504644, 0.02299698679504644, 1.300563710985, 11.705073398865, 0.01386544842799324, 0.01386544842799324)
  ,('FAPs', 'Efnb2', 'Ephb1', 'ECs', 3, 1, 11.420477, 34.261431, 0.1838793176915316, 0.1838793176915316, 3, 1, 1.397474666666667, 4.192424, 0.9253320219903044, 0.9253320219903044, 15.95982728874933, 143.638445598744, 0.1701494208417025, 0.1701494208417025)
  ,('FAPs', 'Efnb2', 'Ephb1', 'MuSCs', 3, 1, 11.420477, 34.261431, 0.1838793176915316, 0.1838793176915316, 3, 1, 0.07803566666666666, 0.234107, 0.05167099121464913, 0.05167099121464913, 0.8912045363463332, 8.020840827117, 0.009501226608994807, 0.009501226608994807)
  ,('FAPs', 'Efnb2', 'Ephb1', 'Neutrophils', 3, 1, 11.420477, 34.261431, 0.1838793176915316, 0.1838793176915316, 1, 0.3333333333333333, 0.034731, 0.104193, 0.02299698679504644, 0.02299698679504644, 0.396644586687, 3.569801280183, 0.004228670240834302, 0.004228670240834302)
  ,('Inflammatory-Mac', 'Efnb2', 'Ephb1', 'ECs', 3, 1, 1.244612333333333, 3.733837, 0.02003930892236799, 0.02003930892236799, 3, 1, 1.397474666666667, 4.192424, 0.9253320219903044, 0.9253320219903044, 1.739314205654222, 15.653827850888, 0.01854301424442312, 0.01854301424442312)
  ,('Inflammatory-Mac', 'Efnb2', 'Ephb1', 'MuSCs', 3, 1, 1.244612333333333, 3.733837, 0.02003930892236799, 0.02003930892236799, 3, 1, 0.07803566666666666, 0.234107, 0.05167099121464913, 0.05167099121464913, 0.0971241531732222, 0.8741173785589998, 0.001035450955275317, 0.001035450955275317)
  ,('Inflammatory-Mac', 'Efnb2', 'Ephb1', 'Neutrophils', 3, 1, 1.244612333333333, 3.733837, 0.02003930892236799, 0.02003930892236799, 1, 0.3333333333333333, 0.034731, 0.104193, 0.02299698679504644, 0.02299698679504644, 0.04322663094899999, 0.389039678541, 0.0004608437226695531, 0.000460843722669553)
  ,('MuSCs', 'Efnb2', 'Ephb1', 'ECs', 3, 1, 10.03858766666667, 30.115763, 0.1616297332180864, 0.1616297332180864, 3, 1, 1.397474666666667, 4.192424, 0.9253320219903044, 0.9253320219903044, 14.02867195327911, 126.258047579512, 0.1495611678524453, 0.1495611678524453)
  ,('MuSCs', 'Efnb2', 'Ephb1', 'MuSCs', 3, 1, 10.03858766666667, 30.115763, 0.1616297332180864, 0.1616297332180864, 3, 1, 0.07803566666666666, 0.234107, 0.05167099121464913, 0.05167099121464913, 0.783367880960111, 7.050310928640999, 0.008351568525137822, 0.008351568525137822)
  ,('MuSCs', 'Efnb2', 'Ephb1', 'Neutrophils', 3, 1, 10.03858766666667, 30.115763, 0.1616297332180864, 0.1616297332180864, 1, 0.3333333333333333, 0.034731, 0.104193, 0.02299698679504644, 0.02299698679504644, 0.348650188251, 3.137851694259, 0.003716996840503211, 0.003716996840503211)
  ,('Neutrophils', 'Efnb2', 'Ephb1', 'ECs', 3, 1, 1.327177333333333, 3.981532, 0.02136867510078605, 0.02136867510078605, 3, 1, 1.397474666666667, 4.192424, 0.9253320219903044, 0.9253320219903044, 1.854696701507556, 16.692270313568, 0.01977311933826423, 0.01977311933826423)
  ,('Neutrophils', 'Efnb2', 'Ephb1', 'MuSCs', 3, 1, 1.327177333333333, 3.981532, 0.02136867510078605, 0.02136867510078605, 3, 1, 0.07803566666666666, 0.234107, 0.05167099121464913, 0.05167099121464913, 0.1035671679915555, 0.9321045119239999, 0.001104140623401408, 0.001104140623401408)
  ,('Neutrophils', 'Efnb2', 'Ephb1', 'Neutrophils', 3, 1, 1.327177333333333, 3.981532, 0.02136867510078605, 0.02136867510078605, 1, 0.3333333333333333, 0.034731, 0.104193, 0.02299698679504644, 0.02299698679504644, 0.046094195964, 0.414847763676, 0.0004914151391204146, 0.0004914151391204144)
  ,('Resolving-Mac', 'Efnb2', 'Ephb1', 'ECs', 3, 1, 0.630923, 1.892769, 0.01015839274978569, 0.01015839274978569, 3, 1, 1.397474666666667, 4.192424, 0.9253320219903044, 0.9253320219903044, 0.8816989091173333, 7.935290182056, 0.00939988610333084, 0.009399886103330838)
  ,('Resolving-Mac', 'Efnb2', 'Ephb1', 'MuSCs', 3, 1, 0.630923, 1.892769, 0.01015839274978569, 0.01015839274978569, 3, 1, 0.07803566666666666, 0.234107, 0.05167099121464913, 0.05167099121464913, 0.04923449692033333, 0.4431104722829999, 0.0005248942225291318, 0.0005248942225291317)
  ,('Resolving-Mac', 'Efnb2', 'Ephb1', 'Neutrophils', 3, 1, 0.630923, 1.892769, 0.01015839274978569, 0.01015839274978569, 1, 0.3333333333333333, 0.034731, 0.104193, 0.02299698679504644, 0.02299698679504644, 0.021912586713, 0.197213280417, 0.000233612423925717, 0.0002336124239257169)
)

$r = 2
foreach ($row in $rows) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c++
  }
  $r++
}
